$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B/C columns for rows 44-51 (new coin inserted, list shifted down) ---
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

# --- Update D column (Price) ---
$ws.Range("D2").Value = "29.392.82"
$ws.Range("D3").Value = "1.848.53"
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D5").Value = "'240.29"
$ws.Range("D6").Value = "'0.6297"
$ws.Range("D8").Value = "'0.07615"
$ws.Range("D9").Value = "'0.2935"
$ws.Range("D10").Value = "'24.51"
$ws.Range("D11").Value = "'0.07743"
$ws.Range("D12").Value = "1.845.38"
$ws.Range("D13").Value = "'5.007"
$ws.Range("D14").Value = "'0.00001081"
$ws.Range("D15").Value = "'0.6789"
$ws.Range("D16").Value = "'83.65"
$ws.Range("D17").Value = "2.093.85"
$ws.Range("D18").Value = "'6.155"
$ws.Range("D19").Value = "29.416.85"
$ws.Range("D20").Value = "'228.89"
$ws.Range("D23").Value = "'7.449"
$ws.Range("D25").Value = "'157.33"
$ws.Range("D26").Value = "'0.1391"
$ws.Range("D27").Value = "'8.379"
$ws.Range("D28").Value = "'17.64"
$ws.Range("D29").Value = "'1.313"
$ws.Range("D31").Value = "'0.05613"
$ws.Range("D33").Value = "'4.047"
$ws.Range("D34").Value = "'1.849"
$ws.Range("D36").Value = "'0.7089"
$ws.Range("D37").Value = "'2.584"
$ws.Range("D38").Value = "1.233.24"
$ws.Range("D39").Value = "'2.773"
$ws.Range("D40").Value = "'0.01799"
$ws.Range("D42").Value = "'0.9093"
$ws.Range("D43").Value = "'0.9999"
$ws.Range("D44").Value = "2.003.08"
$ws.Range("D45").Value = "'101.56"
$ws.Range("D46").Value = "'66.03"
$ws.Range("D47").Value = "'0.00000000121"
$ws.Range("D48").Value = "'7.159"
$ws.Range("D49").Value = "'0.4017"
$ws.Range("D50").Value = "'8.986"
$ws.Range("D51").Value = "'1.684"

# --- Update E column (Volume 1h) ---
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  +7.17%  "
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("E17").Value = "  -7.53%  "
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +4.96%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("E41").Value = "  +4.33%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("E47").Value = "  +3.15%  "
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("E51").Value = "  -1.44%  "
